# DCO_out_freq_calculator.xlsx update
# - Insert a new "Relative Error [%]" row between the real_out_freq row and
#   the CLKSRC row.
# - Re-label in_freq / out_freq / real_out_freq rows to show units ("[Hz]").
# - Change the desired output-frequency (resolution) target from 25000 to 1000.
# - Re-style the "out_freq [Hz]" input cell (C5) as a bold highlighted input box
#   and re-style the CLKSRC row to match the other computed (green) rows.
# - Update the saved selection to J18 (was J12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 8 -- this pushes the old row 8 (CLKSRC) down to row 9,
#    old row 9 (blank) down to row 10, old row 13/14 down to 14/15, and carries
#    every formula/relative-reference along with it.
$ws.Rows("8:8").Insert()

# 2) Re-label the rows that now show units. (Order matches the shared-string
#    table layout of the authored workbook: real_out_freq, out_freq, in_freq.)
$ws.Range("B7").Value = "real_out_freq [Hz]"
$ws.Range("B5").Value = "out_freq [Hz]"
$ws.Range("B4").Value = "in_freq [Hz]"

# 3) New row 8: "Relative Error [%]" label + formulas.
$ws.Range("B8").Value = "Relative Error [%]"
$ws.Range("C8").Formula = "=((C7-`$C`$5)/`$C`$5)*100"
$ws.Range("D8").Formula = "=((D7-`$C`$5)/`$C`$5)*100"
$ws.Range("E8").Formula = "=((E7-`$C`$5)/`$C`$5)*100"
$ws.Range("F8").Formula = "=((F7-`$C`$5)/`$C`$5)*100"
$ws.Range("G8").Formula = "=((G7-`$C`$5)/`$C`$5)*100"

# 4) Restyle B8:G8 (the new Relative Error row) like the B4:G4 row above it:
#    same label style as the other B-column row headers, same fill/border as
#    the numeric in_freq row for the value cells. (Formats only -- formulas
#    already in place are left untouched.)
$ws.Range("B4:G4").Copy() | Out-Null
$ws.Range("B8:G8").PasteSpecial(-4122) | Out-Null

# Put the bold row-label style back (B8 keeps the bold "B-column" look).
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null

# 5) Restyle C5 (out_freq target) as a highlighted bold input box matching the
#    workbook's header box look, then set its value (formats-only paste keeps
#    the 1000 write below intact).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null

# 6) Change the target output frequency / resolution value.
$ws.Range("C5").Value = 1000

# 7) Restyle the CLKSRC row (now row 9) like the other computed/green rows
#    (e.g. CNT_VAL row 6), keeping its own bottom border + center alignment.
$ws.Range("C6:G6").Copy() | Out-Null
$ws.Range("C9:G9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9:G9").HorizontalAlignment = -4108
$ws.Range("C9:G9").Borders.Item(9).LineStyle = 1
$ws.Range("C9:G9").Borders.Item(9).Weight = -4138

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null

# 8) Restore formula text after the paste-special formatting passes
#    (PasteSpecial(xlPasteFormats) should not disturb formulas, but make sure
#    the CLKSRC formulas are untouched).
$ws.Range("C9").Formula = "=DEC2BIN(0,3)"
$ws.Range("D9").Formula = "=DEC2BIN(1,3)"
$ws.Range("E9").Formula = "=DEC2BIN(2,3)"
$ws.Range("F9").Formula = "=DEC2BIN(3,3)"
$ws.Range("G9").Formula = "=DEC2BIN(4,3)"

# 9) Restore the saved selection.
$ws.Range("J18").Select() | Out-Null

$excel.CutCopyMode = 0
